$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Priority column (E) for the four files that just got handed off: low -> ht
foreach ($r in 4..7) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}

# Latest Handoff Datetime (H) for zh-cn rows 4-7: bump the handoff timestamp
foreach ($r in 4..7) {
    $wsZhCn.Range("H$r").Value = "2016-09-01 22:33:37"
}

# Latest HO Xliff Generate Date / de-de Latest Handoff Datetime: bump it too.
# This text is shared between the Overview sheet (col G) and the de-de sheet
# (col H) for the same four rows, so update every cell that references it.
foreach ($r in 4..7) {
    $wsOverview.Range("G$r").Value = "2016-09-01 22:33:42"
    $wsDeDe.Range("H$r").Value = "2016-09-01 22:33:42"
}
